# Add a new reference row (row 8) to the prey ingest / predict table,
# citing the authors' own estimate ("Savoca et al., this study").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 0.75
$ws.Range("D8").Value = "Savoca et al., this study"

# Move/update the active selection to reflect the next empty row, as in
# the authored workbook.
$ws.Range("B9").Select()
